$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# Row 2 used to hold the subtitle "(according to the population census
# data)". The new layout drops that text but keeps the (now empty) row.
$ws.Range("A2").ClearContents()

# Row 3 was a blank spacer row that separated the title block from the
# "(sq. km)" label; the new layout removes that spacer so the label moves
# up to row 3.
$ws.Rows.Item(3).Delete()

# The data table used to show three census years (1989, 2002, 2014); the
# new layout keeps only the most recent year (2014), so drop columns B
# and C - column D (2014 figures) shifts left to become column B.
$ws.Range("B:C").Delete()

# --- Row heights -------------------------------------------------------
# Every row in the refreshed layout uses an explicit 20.1pt row height.
for ($r = 1; $r -le 8; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.1
}

# --- Make sure the trailing blank rows are part of the sheet's used
# range (A1:B8) the same way the refreshed template does, without
# leaving stray values behind.
$ws.Range("A6:B8").Formula = '=""'
$ws.Range("A6:B8").ClearContents()

Write-Output "Shuakhevi municipality area sheet updated"
